$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Fill in the three new agenda rows (11-13) that were previously blank.
$ws.Range("A11").Value = "Pedro"
$ws.Range("B11").Value = "'2284"
$ws.Range("C11").Value = "Igreja Betel"
$ws.Range("D11").Value = "Zona aberta"
$ws.Range("G11").Value = "Pendente"

$ws.Range("A12").Value = "Pedro"
$ws.Range("B12").Value = "'2194"
$ws.Range("C12").Value = "Anselmo"
$ws.Range("D12").Value = "Cliente pedindo revisão no sistema de alarme."
$ws.Range("G12").Value = "Pendente"

$ws.Range("A13").Value = "Pedro"
$ws.Range("B13").Value = "'2029"
$ws.Range("C13").Value = "Marco Otavio"
$ws.Range("D13").Value = "Comunicação instável, cliente pedindo revisão."
$ws.Range("G13").Value = "Pendente"

# Update the view / selection state to match the saved workbook.
$ws.Activate()
$appWindow = $excel.ActiveWindow
$appWindow.ScrollColumn = 6
$ws.Range("H13").Select()
